# Applies scheduled market-data refresh values to the Leve profit sheets.
# Source: xml diff against Sheets/Belias_Profits.xlsx (per-sheet: ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 444.66666
$ws.Range("I4").Value = 376.42856
$ws.Range("K4").Value = 376.42856
$ws.Range("M4").Value = -262.42856

$ws.Range("H8").Value = 277.83334
$ws.Range("I8").Value = 167.5
$ws.Range("J8").Value = 498.5
$ws.Range("K8").Value = 502.5
$ws.Range("L8").Value = 1495.5
$ws.Range("M8").Value = -363.5
$ws.Range("N8").Value = -1773.5

$ws.Range("H18").Value = 639.6
$ws.Range("I18").Value = 733
$ws.Range("J18").Value = 499.5
$ws.Range("K18").Value = 733
$ws.Range("L18").Value = 499.5
$ws.Range("M18").Value = -449
$ws.Range("N18").Value = -1067.5

$ws.Range("H28").Value = 1715.6923
$ws.Range("I28").Value = 1301.2632
$ws.Range("K28").Value = 1301.2632
$ws.Range("M28").Value = -816.2632000000001

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H112").Value = 1112.625
$ws.Range("J112").Value = 1171.5714
$ws.Range("L112").Value = 3514.7142
$ws.Range("N112").Value = -5730.7142

$ws.Range("H113").Value = 4511.3794
$ws.Range("I113").Value = 4432.5
$ws.Range("J113").Value = 4686.6665
$ws.Range("K113").Value = 4432.5
$ws.Range("L113").Value = 4686.6665
$ws.Range("M113").Value = -1178.5
$ws.Range("N113").Value = -11194.6665

$ws.Range("H135").Value = 35715400
$ws.Range("I135").Value = 859.3182
$ws.Range("J135").Value = 166668720
$ws.Range("K135").Value = 7733.8638
$ws.Range("L135").Value = 1500018480
$ws.Range("M135").Value = -5198.8638
$ws.Range("N135").Value = -1500023550

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4330695
$ws.Range("I45").Value = 4786368
$ws.Range("J45").Value = 1800
$ws.Range("K45").Value = 4786368
$ws.Range("L45").Value = 1800
$ws.Range("M45").Value = -4785991
$ws.Range("N45").Value = -2554

$ws.Range("H118").Value = 29928.572
$ws.Range("J118").Value = 29928.572
$ws.Range("L118").Value = 29928.572
$ws.Range("N118").Value = -33242.572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 30000
$ws.Range("J50").Value = 30000
$ws.Range("L50").Value = 30000
$ws.Range("N50").Value = -31148

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 660
$ws.Range("I16").Value = 575
$ws.Range("K16").Value = 575
$ws.Range("M16").Value = -288

$ws.Range("H96").Value = 28200
$ws.Range("J96").Value = 28200
$ws.Range("L96").Value = 28200
$ws.Range("N96").Value = -33692

$ws.Range("H113").Value = 660
$ws.Range("I113").Value = 575
$ws.Range("K113").Value = 575
$ws.Range("M113").Value = 1595

$ws.Range("H132").Value = 927946.1
$ws.Range("I132").Value = 1824.875
$ws.Range("J132").Value = 4632431
$ws.Range("K132").Value = 5474.625
$ws.Range("L132").Value = 13897293
$ws.Range("M132").Value = -2944.625
$ws.Range("N132").Value = -13902353

$ws.Range("H134").Value = 3336.348
$ws.Range("I134").Value = 2823
$ws.Range("J134").Value = 5774.75
$ws.Range("K134").Value = 8469
$ws.Range("L134").Value = 17324.25
$ws.Range("M134").Value = -5934
$ws.Range("N134").Value = -22394.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 697.0263
$ws.Range("I5").Value = 312.30435
$ws.Range("K5").Value = 936.91305
$ws.Range("M5").Value = -824.91305

$ws.Range("H131").Value = 923.4400000000001
$ws.Range("J131").Value = 954.7234
$ws.Range("L131").Value = 2864.1702
$ws.Range("N131").Value = -12944.1702

$ws.Range("H132").Value = 963483.6
$ws.Range("I132").Value = 2040
$ws.Range("J132").Value = 3126731.8
$ws.Range("K132").Value = 18360
$ws.Range("L132").Value = 28140586.2
$ws.Range("M132").Value = -15830
$ws.Range("N132").Value = -28145646.2

$ws.Range("H135").Value = 697.0263
$ws.Range("I135").Value = 312.30435
$ws.Range("K135").Value = 2810.73915
$ws.Range("M135").Value = -275.7391499999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 1330
$ws.Range("J9").Value = 3000
$ws.Range("L9").Value = 3000
$ws.Range("N9").Value = -3340

$ws.Range("H92").Value = 8917
$ws.Range("J92").Value = 8917
$ws.Range("L92").Value = 8917
$ws.Range("N92").Value = -12661

$ws.Range("H102").Value = 37038704
$ws.Range("I102").Value = 111111110
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 111111110
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -111109488
$ws.Range("N102").Value = -5744

$ws.Range("H113").Value = 1900
$ws.Range("I113").Value = 1900
$ws.Range("K113").Value = 1900
$ws.Range("M113").Value = 270

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1801.3334
$ws.Range("I7").Value = 1802
$ws.Range("K7").Value = 1802
$ws.Range("M7").Value = -1690

$ws.Range("H40").Value = 2991.389
$ws.Range("I40").Value = 2865.3125
$ws.Range("K40").Value = 2865.3125
$ws.Range("M40").Value = -2729.3125

$ws.Range("H126").Value = 1801.3334
$ws.Range("I126").Value = 1802
$ws.Range("K126").Value = 5406
$ws.Range("M126").Value = -2936

$ws.Range("H127").Value = 32183.166
$ws.Range("J127").Value = 32183.166
$ws.Range("L127").Value = 32183.166
$ws.Range("N127").Value = -42103.166

$ws.Range("H132").Value = 3611.1304
$ws.Range("I132").Value = 3685.3704
$ws.Range("J132").Value = 3505.6316
$ws.Range("K132").Value = 11056.1112
$ws.Range("L132").Value = 10516.8948
$ws.Range("M132").Value = -8526.111199999999
$ws.Range("N132").Value = -15576.8948

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 989.5714
$ws.Range("I126").Value = 963
$ws.Range("J126").Value = 1042.7142
$ws.Range("K126").Value = 2889
$ws.Range("L126").Value = 3128.1426
$ws.Range("M126").Value = -419
$ws.Range("N126").Value = -8068.142599999999

$ws.Range("H138").Value = 31571.75
$ws.Range("J138").Value = 31571.75
$ws.Range("L138").Value = 31571.75
$ws.Range("N138").Value = -41851.75
